# This edit reorders the data rows (rows 2-76, columns A-K) of the "maxaro"
# product sheet. The header row (row 1) is untouched; the set of data rows
# is identical (same 75 products), just rearranged into a new order.
#
# $order[i-1] gives the 1-based position (within the ORIGINAL A2:K76 block)
# of the row whose contents should become the i-th row of the NEW A2:K76
# block (i.e. new worksheet row (i+1)).
$order = @(6,1,3,10,2,7,13,8,5,12,4,11,9,14,19,18,15,17,16,21,20,24,22,25,26,28,23,27,30,33,31,29,32,39,34,41,35,36,37,40,38,42,44,47,43,45,51,48,46,50,57,53,52,54,55,49,58,59,60,61,64,65,56,62,66,63,67,69,68,71,70,73,72,75,74)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 76
$firstCol = 1   # A
$lastCol = 11   # K

$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$data = $srcRange.Value2

$rowCount = $data.GetUpperBound(0)
$colCount = $data.GetUpperBound(1)

# Snapshot the original rows as plain PowerShell arrays (1 per worksheet row)
# so that we can freely overwrite $data afterwards without losing source
# values (the mapping is a permutation, so naive in-place swapping would
# clobber values before they are read).
$rows = @()
for ($r = 1; $r -le $rowCount; $r++) {
    $rowArr = @()
    for ($c = 1; $c -le $colCount; $c++) {
        $rowArr += $data[$r, $c]
    }
    $rows += ,$rowArr
}

for ($r = 1; $r -le $rowCount; $r++) {
    $srcIdx = $order[$r - 1] - 1
    $srcRow = $rows[$srcIdx]
    for ($c = 1; $c -le $colCount; $c++) {
        $data[$r, $c] = $srcRow[$c - 1]
    }
}

$srcRange.Value2 = $data
